# Tracking-list workbook update: append 12 new apartment-complex rows
# (rows 241-252) to Sheet1, matching the styling already used by the
# existing rows at the bottom of the list, then move the viewport /
# selection to reflect where the user ended up after typing the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Clone existing cell formatting onto the new rows before we put
#    any values in, so the new cells pick up the same cellXfs (font,
#    fill, border, alignment) as their "siblings" higher up the sheet.
# ---------------------------------------------------------------

# Column A, rows 241-248: same look as A239 (right-aligned, no fill,
# boxed border).
$ws.Range("A239").Copy()
$ws.Range("A241:A248").PasteSpecial(-4122)

# Column A, rows 249, 251, 252: same look as A217 (centered, shaded
# fill, boxed border).
$ws.Range("A217").Copy()
$ws.Range("A249").PasteSpecial(-4122)
$ws.Range("A251").PasteSpecial(-4122)
$ws.Range("A252").PasteSpecial(-4122)

# Column A, row 250: same look as A218 (centered, white fill, boxed
# border).
$ws.Range("A218").Copy()
$ws.Range("A250").PasteSpecial(-4122)

# Column B, rows 241-252: new look - same font/border as A239's style
# but without the right alignment (defaults to general).
$ws.Range("A239").Copy()
$ws.Range("B241:B252").PasteSpecial(-4122)
$ws.Range("B241:B252").HorizontalAlignment = 1

$excel.CutCopyMode = 0

# Match the row height used by the rest of the formatted table.
$ws.Range("A241:B252").RowHeight = 18

# ---------------------------------------------------------------
# 2. Write the new rows of data (id in column A, apartment name in
#    column B).
# ---------------------------------------------------------------
$ws.Range("A241").Value2 = 3615
$ws.Range("B241").Value2 = "진달래효성"

$ws.Range("A242").Value2 = 8247
$ws.Range("B242").Value2 = "다정한마을KCC"

$ws.Range("A243").Value2 = 1420
$ws.Range("B243").Value2 = "미리내동성"

$ws.Range("A244").Value2 = 103865
$ws.Range("B244").Value2 = "부천소사푸르지오"

$ws.Range("A245").Value2 = 14462
$ws.Range("B245").Value2 = "소새울역중흥S클래스"

$ws.Range("A246").Value2 = 110681
$ws.Range("B246").Value2 = "부천옥길호반베르디움"

$ws.Range("A247").Value2 = 102622
$ws.Range("B247").Value2 = "역곡역e편한세상"

$ws.Range("A248").Value2 = 127082
$ws.Range("B248").Value2 = "일루미스테이트"

$ws.Range("A249").Value2 = 27540
$ws.Range("B249").Value2 = "두산위브트레지움2단지"

$ws.Range("A250").Value2 = 120265
$ws.Range("B250").Value2 = "e편한세상온수"

$ws.Range("A251").Value2 = 108756
$ws.Range("B251").Value2 = "송내역파인푸르지오1단지"

$ws.Range("A252").Value2 = 102634
$ws.Range("B252").Value2 = "원종금호어울림"

# ---------------------------------------------------------------
# 3. Reflect the final scroll/selection state: the user scrolled down
#    so row 246 is at the top and landed just past the new data.
# ---------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 246
$ws.Range("E256").Select() | Out-Null
